$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: B1 should get the value currently in D1 (Alsto), then D1 is cleared.
$ws.Range("B1").Value = $ws.Range("D1").Value2

# Row 21: B21 should get the value currently in D21 (Aries), then D21 is cleared.
$ws.Range("B21").Value = $ws.Range("D21").Value2

# For every other data row (2-20, 22-25, 27-33), column B gets the same
# value as column A in that row.
$rows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,22,23,24,25,27,28,29,30,31,32,33)
foreach ($r in $rows) {
    $aAddr = "A" + $r
    $bAddr = "B" + $r
    $ws.Range($bAddr).Value = $ws.Range($aAddr).Value2
}

# Row 26 already has B26 populated (Rose) and D26 duplicates it; just clear D26.
$ws.Range("D26").ClearContents()

# Column D is no longer used anywhere; clear out the now-stale D1/D21 cells
# (values were already relocated to B above).
$ws.Range("D1").ClearContents()
$ws.Range("D21").ClearContents()

# Add a brand-new row 34 with a new actor name in column A only.
$ws.Range("A34").Value = "リリー"
